# Update cryptocurrency price/volume data per GitHub Actions refresh
# (values assigned as leading-apostrophe strings are numeric-looking
#  price text, e.g. "210.99" - the apostrophe forces Excel to keep
#  them as text instead of auto-converting to a Number)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.135.35'
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").Value = '1.669.86'
$ws.Range("E3").Value = '  -1.44%  '

$ws.Range("E4").Value = '  -0.78%  '

$ws.Range("D5").Value = "'210.99"
$ws.Range("E5").Value = '  -3.80%  '

$ws.Range("D6").Value = "'0.5225"
$ws.Range("E6").Value = '  -4.64%  '

$ws.Range("E7").Value = '  -0.76%  '

$ws.Range("D8").Value = "'0.2642"
$ws.Range("E8").Value = '  -3.43%  '

$ws.Range("D9").Value = "'0.06261"
$ws.Range("E9").Value = '  -3.22%  '

$ws.Range("D10").Value = "'21.16"
$ws.Range("E10").Value = '  -3.88%  '

$ws.Range("D11").Value = "'0.07517"
$ws.Range("E11").Value = '  -2.10%  '

$ws.Range("D12").Value = '1.685.43'
$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").Value = "'4.440"
$ws.Range("E13").Value = '  -2.54%  '

$ws.Range("D14").Value = "'0.5596"
$ws.Range("E14").Value = '  -4.37%  '

$ws.Range("D15").Value = "'0.000007974"
$ws.Range("E15").Value = '  -4.77%  '

$ws.Range("D16").Value = "'66.31"
$ws.Range("E16").Value = '  +1.27%  '

$ws.Range("D17").Value = '26.185.87'
$ws.Range("E17").Value = '  -0.98%  '

$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("D19").Value = "'4.787"
$ws.Range("E19").Value = '  -3.27%  '

$ws.Range("D20").Value = "'187.11"
$ws.Range("E20").Value = '  -2.84%  '

$ws.Range("E21").Value = '  -5.64%  '

$ws.Range("D22").Value = "'6.175"
$ws.Range("E22").Value = '  -1.32%  '

$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = '  -0.74%  '

$ws.Range("D24").Value = "'147.90"
$ws.Range("E24").Value = '  -0.95%  '

$ws.Range("D25").Value = "'0.1245"
$ws.Range("E25").Value = '  -6.34%  '

$ws.Range("D26").Value = "'7.587"
$ws.Range("E26").Value = '  -4.15%  '

$ws.Range("D27").Value = "'15.92"
$ws.Range("E27").Value = '  +0.97%  '

$ws.Range("D28").Value = "'0.06204"
$ws.Range("E28").Value = '  -1.41%  '

$ws.Range("D29").Value = "'1.356"
$ws.Range("E29").Value = '  -2.84%  '

$ws.Range("D30").Value = "'1.279"
$ws.Range("E30").Value = '  -3.84%  '

$ws.Range("D31").Value = "'3.471"
$ws.Range("E31").Value = '  -3.97%  '

$ws.Range("D32").Value = "'3.429"
$ws.Range("E32").Value = '  -4.89%  '

$ws.Range("D33").Value = "'1.611"
$ws.Range("E33").Value = '  -4.41%  '

$ws.Range("D34").Value = "'0.9918"
$ws.Range("E34").Value = '  -5.07%  '

$ws.Range("D35").Value = "'0.6052"
$ws.Range("E35").Value = '  -1.49%  '

$ws.Range("D36").Value = "'2.403"
$ws.Range("E36").Value = '  -0.32%  '

$ws.Range("D37").Value = "'2.708"
$ws.Range("E37").Value = '  -0.09%  '

$ws.Range("D38").Value = "'6.130"
$ws.Range("E38").Value = '  -1.09%  '

$ws.Range("D39").Value = "'0.01611"
$ws.Range("E39").Value = '  -1.74%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'0.8657"
$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.068.99'
$ws.Range("E41").Value = '  -4.49%  '

$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = '  -1.13%  '

$ws.Range("D43").Value = "'99.58"
$ws.Range("E43").Value = '  -2.31%  '

$ws.Range("D44").Value = '1.819.06'
$ws.Range("E44").Value = '  -1.38%  '

$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("E45").Value = '  -0.95%  '

$ws.Range("D46").Value = "'55.95"
$ws.Range("E46").Value = '  -2.75%  '

$ws.Range("E47").Value = '  -0.32%  '

$ws.Range("D48").Value = "'0.05248"
$ws.Range("E48").Value = '  -0.69%  '

$ws.Range("D49").Value = "'7.923"
$ws.Range("E49").Value = '  -3.45%  '

$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("D51").Value = "'5.950"
$ws.Range("E51").Value = '  -2.51%  '
